$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the next empty row after the existing data (row 74 -> new row 75)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

$ws.Cells.Item($newRow, 1).Value = "Glycolysis and Gluconeogenesis"
$ws.Cells.Item($newRow, 2).Value = "Carbohydrate metabolism"

# Update the view/selection similar to what the author ended up with when saving
$ws.Range("A26").Select() | Out-Null
$ws.Range("A76").Select() | Out-Null
